# Update gh-pages output (F column counters refreshed, plus the F2 date
# on sheet1 "展览") to the values generated at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# SHEET1
$ws1.Range("F2").Value = 42315
$ws1.Range("F4").Value = 10044
$ws1.Range("F5").Value = 227
$ws1.Range("F6").Value = 1063
$ws1.Range("F7").Value = 982
$ws1.Range("F8").Value = 781
$ws1.Range("F9").Value = 243
$ws1.Range("F10").Value = 317
$ws1.Range("F11").Value = 1023
$ws1.Range("F13").Value = 138
$ws1.Range("F15").Value = 350
$ws1.Range("F16").Value = 1627
$ws1.Range("F17").Value = 15
$ws1.Range("F18").Value = 810
$ws1.Range("F19").Value = 761
$ws1.Range("F20").Value = 504
$ws1.Range("F21").Value = 726
$ws1.Range("F22").Value = 806
$ws1.Range("F24").Value = 259
$ws1.Range("F25").Value = 74
$ws1.Range("F26").Value = 576
$ws1.Range("F28").Value = 80
$ws1.Range("F29").Value = 279
$ws1.Range("F30").Value = 968
$ws1.Range("F32").Value = 460
$ws1.Range("F33").Value = 123
$ws1.Range("F35").Value = 177
$ws1.Range("F36").Value = 487
$ws1.Range("F37").Value = 1438
$ws1.Range("F38").Value = 326
$ws1.Range("F39").Value = 1318
$ws1.Range("F40").Value = 389
$ws1.Range("F42").Value = 31
$ws1.Range("F43").Value = 53
$ws1.Range("F45").Value = 54
$ws1.Range("F46").Value = 18
$ws1.Range("F47").Value = 19
$ws1.Range("F48").Value = 4

# SHEET2
$ws2.Range("F3").Value = 351
$ws2.Range("F4").Value = 4482
$ws2.Range("F8").Value = 159
$ws2.Range("F16").Value = 45
$ws2.Range("F17").Value = 4392

# SHEET3
$ws3.Range("F2").Value = 2112
$ws3.Range("F3").Value = 572
$ws3.Range("F4").Value = 477

# SHEET4
$ws4.Range("F2").Value = 2112
$ws4.Range("F3").Value = 572
$ws4.Range("F4").Value = 351
$ws4.Range("F6").Value = 10044
$ws4.Range("F7").Value = 1063
$ws4.Range("F8").Value = 1063
$ws4.Range("F10").Value = 477
$ws4.Range("F11").Value = 982
$ws4.Range("F12").Value = 781
$ws4.Range("F13").Value = 159
$ws4.Range("F14").Value = 317
$ws4.Range("F15").Value = 1023
$ws4.Range("F18").Value = 350
$ws4.Range("F19").Value = 1627
$ws4.Range("F20").Value = 15
$ws4.Range("F21").Value = 810
$ws4.Range("F22").Value = 761
$ws4.Range("F23").Value = 504
$ws4.Range("F24").Value = 726
$ws4.Range("F25").Value = 806
$ws4.Range("F27").Value = 259
$ws4.Range("F28").Value = 74
$ws4.Range("F29").Value = 576
$ws4.Range("F32").Value = 80
$ws4.Range("F33").Value = 279
$ws4.Range("F34").Value = 968
$ws4.Range("F37").Value = 460
$ws4.Range("F38").Value = 123
$ws4.Range("F40").Value = 177
$ws4.Range("F42").Value = 326
$ws4.Range("F43").Value = 1318
$ws4.Range("F44").Value = 389
$ws4.Range("F46").Value = 53
$ws4.Range("F49").Value = 19
